$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 474, shifting rows 474:505
# down to 475:506 (so the dimension grows from A1:R505 to A1:R506).
$ws.Rows(474).Insert()

# Populate the newly inserted row 474 with the new data record.
$ws.Cells.Item(474, 1).Value = 4
$ws.Cells.Item(474, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(474, 3).Value = "Los Lagos"
$ws.Cells.Item(474, 4).Value = 45265
$ws.Cells.Item(474, 5).Value = 10
$ws.Cells.Item(474, 6).Value = 100112032
$ws.Cells.Item(474, 7).Value = "Zapallo italiano"
$ws.Cells.Item(474, 8).Value = "Sin especificar"
$ws.Cells.Item(474, 9).Value = "Primera"
$ws.Cells.Item(474, 10).Value = 250
$ws.Cells.Item(474, 11).Value = 18000
$ws.Cells.Item(474, 12).Value = 18000
$ws.Cells.Item(474, 13).Value = 18000
$ws.Cells.Item(474, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(474, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(474, 16).Value = 360
$ws.Cells.Item(474, 17).Value = 50
$ws.Cells.Item(474, 18).Value = "Hortaliza"
